# Insert a new weekly record at row 80 ("Fruta / hortaliza, semanal").
# This pushes the existing rows 80-127 down to 81-128 and keeps the
# worksheet structure (styles, header, etc.) otherwise unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(80).Insert()

$ws.Cells.Item(80, 1).Value = 10
$ws.Cells.Item(80, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(80, 3).Value = "La Araucanía"
$ws.Cells.Item(80, 4).Value = 45126
$ws.Cells.Item(80, 5).Value = 9
$ws.Cells.Item(80, 6).Value = 100112010
$ws.Cells.Item(80, 7).Value = "Achicoria"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 100
$ws.Cells.Item(80, 11).Value = 9000
$ws.Cells.Item(80, 12).Value = 9000
$ws.Cells.Item(80, 13).Value = 9000
$ws.Cells.Item(80, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(80, 15).Value = "Región Metropolitana"
$ws.Cells.Item(80, 16).Value = 500
$ws.Cells.Item(80, 17).Value = 18
$ws.Cells.Item(80, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of the table.
$ws.Cells.Item(80, 4).NumberFormat = $ws.Cells.Item(81, 4).NumberFormat
